# Fruta / hortaliza, semanal
#
# Adds two new weekly price observations (each reported as a "Primera"/
# "Segunda" pair of rows) into the flat Kiwi price log for "Feria
# Lagunitas de Puerto Montt". All pre-existing rows keep their original
# values and are simply pushed down to make room for the four new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param(
        [int]$Row,
        [object[]]$Values
    )
    for ($i = 0; $i -lt $Values.Length; $i++) {
        $ws.Cells.Item($Row, $i + 1).Value = $Values[$i]
    }
    # Column D (Fecha) carries the same date/time number format as the
    # rest of the column.
    $ws.Cells.Item($Row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

# --- Insert the first new pair of rows at row 253 -------------------------
$ws.Rows.Item(253).Resize(2).Insert()

Set-RowValues 253 @(4, "Feria Lagunitas de Puerto Montt", "Los Lagos", 44748, 10, "Fruta", 100101, "Berries", 100101007, "Kiwi", "Hayward", "Primera", 200, 14000, 14000, 14000, "`$/caja 15 kilos", "Región de O'Higgins", 933, 15)
Set-RowValues 254 @(4, "Feria Lagunitas de Puerto Montt", "Los Lagos", 44748, 10, "Fruta", 100101, "Berries", 100101007, "Kiwi", "Hayward", "Segunda", 200, 12000, 12000, 12000, "`$/caja 15 kilos", "Región de O'Higgins", 800, 15)

# --- Insert the second new pair of rows at (current) row 277 -------------
$ws.Rows.Item(277).Resize(2).Insert()

Set-RowValues 277 @(4, "Feria Lagunitas de Puerto Montt", "Los Lagos", 44747, 10, "Fruta", 100101, "Berries", 100101007, "Kiwi", "Hayward", "Primera", 300, 15000, 15000, 15000, "`$/caja 15 kilos", "Región de O'Higgins", 1000, 15)
Set-RowValues 278 @(4, "Feria Lagunitas de Puerto Montt", "Los Lagos", 44747, 10, "Fruta", 100101, "Berries", 100101007, "Kiwi", "Hayward", "Segunda", 300, 12000, 12000, 12000, "`$/caja 15 kilos", "Región de O'Higgins", 800, 15)

Write-Host ("Used range after edits: {0}" -f $ws.UsedRange.Address())
